# Applies the "Tabellendaten aktualisiert für Testzwecke" edit:
# updates a handful of data cells on the (single) worksheet and
# moves the active selection / scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data updates -----------------------------------------------------
$ws.Range("F2").Value = 0.8
$ws.Range("F3").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("F6").Value = 0.65
$ws.Range("F9").Value = 0.62
$ws.Range("C12").Value = 500000
$ws.Range("C13").Value = 500000

# --- view / selection change -------------------------------------------
$ws.Range("D7").Select()
